# daily auto push: 2026-02-12 05:07 UTC
#
# A new reading was logged for 2026/02/12 (木) that belongs right after the
# existing 2026/02/12 rows (808, 809) and before the 2026/12/29 block.
# Insert a new row at 810, shifting everything from the old row 810 down
# through 851 to 811-852, then populate the new row with the logged values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 810:851 down to 811:852, creating a blank row at 810.
$ws.Rows.Item(810).Insert()

# Column A holds the date as literal text (e.g. "2026/12/29"), matching the
# rest of the column. Format the cell as Text first so Excel doesn't
# auto-convert the "2026/02/12" string into a date serial, then restore the
# default (Normal) style so no stray formatting is left on the cell.
$ws.Cells.Item(810, 1).NumberFormat = "@"
$ws.Cells.Item(810, 1).Value = "2026/02/12"
$ws.Cells.Item(810, 1).Style = "Normal"

$ws.Cells.Item(810, 2).Value = "木"
$ws.Cells.Item(810, 3).Value = 13
$ws.Cells.Item(810, 4).Value = 201
